$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.800.05"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.536.50"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "621.81"
$ws.Range("E5").Value = "  +3.97%  "
$ws.Range("D6").Value = "172.44"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "3.533.79"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "7.22"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").Value = "0.586"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "46.28"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "4.095.43"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "8.44"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "607.60"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "3.538.27"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "70.868.66"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "17.72"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "0.881"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "9.12"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").Value = "15.60"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "97.51"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("D29").Value = "33.70"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "8.13"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "3.01"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "6.82"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").Value = "620.79"
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("D36").Value = "0.0495"
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("D37").Value = "10.86"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value = "0.0997"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "56.78"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  -5.64%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "3.347.16"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "0.0₃0725"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("D45").Value = "0.311"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "31.95"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "134.07"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.156"
$ws.Range("E51").Value = "  +7.35%  "
